$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "40.088.45"
Set-TextValue $ws "E2" "  +0.19%  "
Set-TextValue $ws "D3" "2.223.70"
Set-TextValue $ws "E3" "  -0.66%  "
Set-TextValue $ws "E4" "  -0.03%  "
Set-TextValue $ws "D5" "290.78"
Set-TextValue $ws "E5" "  -1.25%  "
Set-TextValue $ws "D6" "88.07"
Set-TextValue $ws "E6" "  +2.30%  "
Set-TextValue $ws "E7" "  -0.40%  "
Set-TextValue $ws "E8" "  -0.10%  "
Set-TextValue $ws "E9" "  +0.67%  "
Set-TextValue $ws "D10" "30.56"
Set-TextValue $ws "E10" "  +0.30%  "
Set-TextValue $ws "D11" "0.0781"
Set-TextValue $ws "E12" "  +3.06%  "
Set-TextValue $ws "D13" "6.50"
Set-TextValue $ws "E13" "  +1.77%  "
Set-TextValue $ws "D14" "2.568.70"
Set-TextValue $ws "E14" "  -0.59%  "
Set-TextValue $ws "E15" "  -1.44%  "
Set-TextValue $ws "D16" "2.204.77"
Set-TextValue $ws "E16" "  -1.31%  "
Set-TextValue $ws "E17" "  +0.81%  "
Set-TextValue $ws "D18" "40.035.18"
Set-TextValue $ws "E18" "  +0.27%  "
Set-TextValue $ws "D19" "11.53"
Set-TextValue $ws "E19" "  +7.62%  "
Set-TextValue $ws "D20" "0.0₃0886"
Set-TextValue $ws "E20" "  -0.73%  "
Set-TextValue $ws "D21" "5.83"
Set-TextValue $ws "E21" "  +0.49%  "
Set-TextValue $ws "E22" "  +0.39%  "
Set-TextValue $ws "D23" "236.13"
Set-TextValue $ws "E23" "  +0.56%  "
Set-TextValue $ws "E24" "  +0.02%  "
Set-TextValue $ws "E25" "  +1.25%  "
Set-TextValue $ws "E26" "  -0.92%  "
Set-TextValue $ws "D27" "22.70"
Set-TextValue $ws "E27" "  -1.46%  "
Set-TextValue $ws "D28" "2.19"
Set-TextValue $ws "E28" "  -0.78%  "
Set-TextValue $ws "E29" "  +0.00%  "
Set-TextValue $ws "D30" "156.25"
Set-TextValue $ws "E30" "  +0.67%  "
Set-TextValue $ws "D31" "31.84"
Set-TextValue $ws "E31" "  -4.99%  "
Set-TextValue $ws "E33" "  +1.93%  "
Set-TextValue $ws "D34" "0.0718"
Set-TextValue $ws "E34" "  +1.00%  "
Set-TextValue $ws "D35" "2.35"
Set-TextValue $ws "E35" "  -0.58%  "
Set-TextValue $ws "E36" "  +6.51%  "
Set-TextValue $ws "E37" "  -0.31%  "
Set-TextValue $ws "D38" "15.80"
Set-TextValue $ws "E38" "  -4.30%  "
Set-TextValue $ws "E39" "  -0.37%  "
Set-TextValue $ws "E40" "  +1.75%  "
Set-TextValue $ws "D41" "2.115.48"
Set-TextValue $ws "E41" "  +8.25%  "
Set-TextValue $ws "D42" "3.84"
Set-TextValue $ws "E42" "  +1.48%  "
Set-TextValue $ws "D43" "2.15"
Set-TextValue $ws "E43" "  -1.66%  "
Set-TextValue $ws "D44" "10.02"
Set-TextValue $ws "E44" "  +5.26%  "
Set-TextValue $ws "B45" "EnergySwap"
Set-TextValue $ws "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D45" "17.93"
Set-TextValue $ws "E45" "  +9.84%  "
Set-TextValue $ws "B46" "VeChain"
Set-TextValue $ws "C46" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D46" "0.0268"
Set-TextValue $ws "E46" "  -1.16%  "
Set-TextValue $ws "E47" "  +2.10%  "
Set-TextValue $ws "D48" "2.434.15"
Set-TextValue $ws "E48" "  -0.87%  "
Set-TextValue $ws "D49" "89.08"
Set-TextValue $ws "E49" "  +0.21%  "
Set-TextValue $ws "B50" "TrustWalletToken"
Set-TextValue $ws "C50" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D50" "1.10"
Set-TextValue $ws "E50" "  +2.76%  "
Set-TextValue $ws "D51" "69.29"
Set-TextValue $ws "E51" "  -2.35%  "
